# Daily refresh of the "剩余" (remaining days) tracker sheet.
# For every data row, recompute the remaining-days count (column E) based on
# the elapsed days since the start date (column F) versus the total days
# (column D). When a cycle completes (remaining days would hit zero or
# below), the start date is rolled forward to "today" and the remaining
# count is reset back to the total. Rows whose start date cannot be parsed
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Today" for this run (one day after the previous refresh).
$today = Get-Date -Year 2025 -Month 11 -Day 1
$todaySerial = [math]::Round($today.ToOADate())

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $totalDays = $dCell.Value2
    $startRaw  = $fCell.Value2

    if ($totalDays -eq $null -or $startRaw -eq $null) {
        continue
    }

    $startStr = [string]([int64]$startRaw)
    if ($startStr.Length -ne 8) {
        # Not a well-formed yyyyMMdd date - skip this row.
        continue
    }

    $year  = [int]$startStr.Substring(0,4)
    $month = [int]$startStr.Substring(4,2)
    $day   = [int]$startStr.Substring(6,2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startDate = Get-Date -Year $year -Month $month -Day $day
    $startSerial = [math]::Round($startDate.ToOADate())

    $elapsed = $todaySerial - $startSerial
    $newRemaining = $totalDays - $elapsed

    if ($newRemaining -le 0) {
        # Cycle finished - restart it as of today.
        $eCell.Value2 = $totalDays
        $fCell.Value2 = ($today.Year * 10000) + ($today.Month * 100) + $today.Day
    } else {
        $eCell.Value2 = $newRemaining
    }
}
